# Update the player_id_y values in column O to reflect the refreshed
# Kaggle CSV data (the IDs shifted slightly after re-downloading).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O2").Value = 3878
$ws.Range("O3").Value = 2716
$ws.Range("O4").Value = 3878
$ws.Range("O5").Value = 3878
$ws.Range("O6").Value = 1951
$ws.Range("O7").Value = 1951
$ws.Range("O8").Value = 2358
$ws.Range("O9").Value = 4457
$ws.Range("O10").Value = 4673
$ws.Range("O11").Value = 4673
$ws.Range("O12").Value = 3094
$ws.Range("O13").Value = 3266
$ws.Range("O14").Value = 3266
$ws.Range("O15").Value = 1346
$ws.Range("O16").Value = 3266
$ws.Range("O17").Value = 3266
$ws.Range("O18").Value = 3154
$ws.Range("O19").Value = 1434
$ws.Range("O20").Value = 4710
$ws.Range("O21").Value = 4710
$ws.Range("O22").Value = 3097
$ws.Range("O23").Value = 4823
$ws.Range("O24").Value = 4823
$ws.Range("O26").Value = 4582
$ws.Range("O27").Value = 2982
$ws.Range("O28").Value = 3649
$ws.Range("O29").Value = 2982
$ws.Range("O30").Value = 3649
$ws.Range("O31").Value = 1241
$ws.Range("O32").Value = 2039
$ws.Range("O33").Value = 804
$ws.Range("O34").Value = 3649
$ws.Range("O35").Value = 3649
$ws.Range("O36").Value = 3409
$ws.Range("O37").Value = 3409
$ws.Range("O38").Value = 3649
$ws.Range("O39").Value = 3409
$ws.Range("O40").Value = 3215
$ws.Range("O41").Value = 3215
$ws.Range("O42").Value = 3215
$ws.Range("O43").Value = 3791
$ws.Range("O44").Value = 3791
$ws.Range("O45").Value = 2937
$ws.Range("O46").Value = 2979
$ws.Range("O47").Value = 3791
$ws.Range("O49").Value = 2979
$ws.Range("O50").Value = 2979
$ws.Range("O52").Value = 2979
$ws.Range("O53").Value = 1200
$ws.Range("O54").Value = 2979
$ws.Range("O55").Value = 2979
$ws.Range("O56").Value = 5226
$ws.Range("O57").Value = 5169
$ws.Range("O58").Value = 5231
$ws.Range("O59").Value = 5231
$ws.Range("O60").Value = 5231
$ws.Range("O62").Value = 3940
$ws.Range("O66").Value = 5231
